$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add Andrew Case's time spent so far
$ws.Range("B2").Value = "8h 30m"
